$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 111541118
$ws.Range("B9").Value = 94851
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 2569
$ws.Range("F9").Value = "Stor revmossa"
$ws.Range("G9").Value = "Bazzania trilobata"
$ws.Range("H9").Value = "(L.) Gray"
$ws.Range("Q9").Value = 693461.6376634488
$ws.Range("R9").Value = 6551559.049034445
$ws.Range("AO9").ClearContents()

# Row 10
$ws.Range("A10").Value = 111541119
$ws.Range("B10").Value = 5426
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 101410
$ws.Range("F10").Value = "Reliktbock"
$ws.Range("G10").Value = "Nothorhina muricata"
$ws.Range("H10").Value = "(Dalman, 1817)"
$ws.Range("M10").Value = "äldre gnagspår"
$ws.Range("Q10").Value = 693467.6220677271
$ws.Range("R10").Value = 6551532.561666255
$ws.Range("AC10").Value = "En gammal tall med kläckhål här och var. Om det är färskt eller gammalt är svårt sia om."
$ws.Range("AO10").Value = "gammeltall"

# Row 11
$ws.Range("A11").Value = 111541120
$ws.Range("B11").Value = 79444
$ws.Range("E11").Value = 1049
$ws.Range("F11").Value = "Kortskaftad ärgspik"
$ws.Range("G11").Value = "Microcalicium ahlneri"
$ws.Range("H11").Value = "Tibell"
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("Q11").Value = 693513.2669972532
$ws.Range("R11").Value = 6551517.868690074
$ws.Range("AF11").ClearContents()
$ws.Range("AO11").Value = "silverstubbe av tall"

# Row 12
$ws.Range("A12").Value = 111541128
$ws.Range("B12").Value = 5113
$ws.Range("E12").Value = 100526
$ws.Range("F12").Value = "Bronshjon"
$ws.Range("G12").Value = "Callidium coriaceum"
$ws.Range("H12").Value = "Paykull, 1800"
$ws.Range("M12").Value = "färska gnagspår"
$ws.Range("Q12").Value = 693570.8046739453
$ws.Range("R12").Value = 6551451.742365629
$ws.Range("AO12").Value = "torrgran"

# Row 14
$ws.Range("A14").Value = 111541129
$ws.Range("B14").Value = 5113
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 100526
$ws.Range("F14").Value = "Bronshjon"
$ws.Range("G14").Value = "Callidium coriaceum"
$ws.Range("H14").Value = "Paykull, 1800"
$ws.Range("M14").Value = "äldre gnagspår"
$ws.Range("Q14").Value = 693328.6441019299
$ws.Range("R14").Value = 6551545.628735202
$ws.Range("AO14").Value = "torrgran"

# Row 15
$ws.Range("A15").Value = 111541121
$ws.Range("B15").Value = 79444
$ws.Range("E15").Value = 1049
$ws.Range("F15").Value = "Kortskaftad ärgspik"
$ws.Range("G15").Value = "Microcalicium ahlneri"
$ws.Range("H15").Value = "Tibell"
$ws.Range("J15").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("Q15").Value = 693460.9606228607
$ws.Range("R15").Value = 6551521.405726598
$ws.Range("AC15").ClearContents()
$ws.Range("AF15").ClearContents()
$ws.Range("AO15").Value = "silverstubbe av tall"

# Row 16
$ws.Range("A16").Value = 111541122
$ws.Range("B16").Value = 5112
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 102204
$ws.Range("F16").Value = "Grönhjon"
$ws.Range("G16").Value = "Callidium aeneum"
$ws.Range("H16").Value = "(De Geer, 1775)"
$ws.Range("Q16").Value = 693344.0451535647
$ws.Range("R16").Value = 6551526.82974836
$ws.Range("AO16").Value = "låga av gran"
